# The sheet is a weekly price log. Two new observation rows need to be
# inserted (pushing the existing rows down by two), with the new rows
# placed right before the current row 69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 69 (shifts old rows 69.. down to 71..)
$ws.Rows.Item(69).EntireRow.Insert()
$ws.Rows.Item(69).EntireRow.Insert()

# New row 69: Ají / Inferno / Primera
$ws.Range("A69").Value = 1
$ws.Range("B69").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C69").Value = "Arica y Parinacota"
$ws.Range("D69").Value = 44806
$ws.Range("E69").Value = 15
$ws.Range("F69").Value = 100112021
$ws.Range("G69").Value = "Ají"
$ws.Range("H69").Value = "Inferno"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 150
$ws.Range("K69").Value = 14000
$ws.Range("L69").Value = 15000
$ws.Range("M69").Value = 14500
$ws.Range("N69").Value = "$/caja 15 kilos"
$ws.Range("O69").Value = "Región de Arica y Parinacota"
$ws.Range("P69").Value = 967
$ws.Range("Q69").Value = 15
$ws.Range("R69").Value = "Hortaliza"

# New row 70: Ají / Inferno / Segunda
$ws.Range("A70").Value = 1
$ws.Range("B70").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C70").Value = "Arica y Parinacota"
$ws.Range("D70").Value = 44806
$ws.Range("E70").Value = 15
$ws.Range("F70").Value = 100112021
$ws.Range("G70").Value = "Ají"
$ws.Range("H70").Value = "Inferno"
$ws.Range("I70").Value = "Segunda"
$ws.Range("J70").Value = 160
$ws.Range("K70").Value = 11000
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = 11500
$ws.Range("N70").Value = "$/caja 15 kilos"
$ws.Range("O70").Value = "Región de Arica y Parinacota"
$ws.Range("P70").Value = 767
$ws.Range("Q70").Value = 15
$ws.Range("R70").Value = "Hortaliza"
